$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Variable" column entry for household income ("A6") was stored as
# " Income " with stray leading/trailing spaces. Replace it with the
# clean "Income" label.
$ws.Range("A6").Value = "Income"

# Reflect the author's resulting selection/view state.
$ws.Range("A6").Select()
$excel.ActiveWindow.Zoom = 120
